$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "checklist"

# Add product codes
$ws.Range("G3").Value = 109
$ws.Range("G6").Value = 108

# Update the active selection to match the diff
$ws.Range("G7").Select()
